# ---------------------------------------------------------------------------
# Commit: Tue, Jul 14, 2020 11:06:43 PM
#
# 1) The table on slide 16 gets a new table style (its <a:tableStyleId> GUID
#    changes from {F3C15FFE-8D38-466F-94E8-0E1979CD64E6} to
#    {BD56528F-7A7D-49B5-981C-5C8C19C0769D}). This is a PowerPoint built-in
#    table style, so it must be applied with Table.ApplyStyle(GUID) rather
#    than by assigning a property.
#
# 2) The presentation's theme palette changes from the "Integral" theme
#    (green/olive accents) back to the default "Office Theme" palette
#    (blue/orange accents). The font scheme and format scheme are already
#    identical between the two theme parts, so only the 12 color slots of
#    the theme color scheme actually need to change.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- helper: "RRGGBB" hex string -> the little-endian long PowerPoint's
#     RGB uses (r + g*256 + b*65536) ------------------------------------
function ComRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# --- 1) Table style on slide 16 -----------------------------------------
$tableSlide = $p.Slides.Item(16)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{BD56528F-7A7D-49B5-981C-5C8C19C0769D}")
    }
}

# --- 2) Theme colors: Integral -> Office Theme --------------------------
# Order matches ThemeColorScheme.Colors(1..12):
#   dk1, lt1, dk2, lt2, accent1..accent6, hlink, folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = ComRGB($officeColors[$i - 1])
}
